$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 0.06366381003880003
$ws.Range("B9").Value = 0
$ws.Range("B16").Value = 0
$ws.Range("B30").Value = 0
$ws.Range("B31").Value = 3.588460746437567
$ws.Range("B34").Value = 0
$ws.Range("B62").Value = 0
$ws.Range("B72").Value = 5
$ws.Range("B73").Value = 0
$ws.Range("B75").Value = 4.154078783986692
$ws.Range("B76").Value = 0
$ws.Range("B78").Value = 4.999999999999999
$ws.Range("B79").Value = 0
$ws.Range("B80").Value = 0
$ws.Range("B81").Value = 4.999999999999998
$ws.Range("B82").Value = 3.243283972901243
$ws.Range("B83").Value = 3.326071396160306
$ws.Range("B84").Value = 5
$ws.Range("B85").Value = 0
$ws.Range("B88").Value = 0
$ws.Range("B91").Value = 0
$ws.Range("B92").Value = 0
$ws.Range("B134").Value = 0
$ws.Range("B153").Value = 0
$ws.Range("B172").Value = 5
$ws.Range("B175").Value = 1.942000742333645
$ws.Range("B176").Value = 3.950490898271326
$ws.Range("B177").Value = 5
$ws.Range("B179").Value = 4.999999999999999
$ws.Range("B181").Value = 4.145487026351776
$ws.Range("B182").Value = 0.2922529296163545
$ws.Range("B183").Value = 5
$ws.Range("B266").Value = 4.596364078714711
$ws.Range("B267").Value = 0
$ws.Range("B268").Value = 4.999999999999998
$ws.Range("B298").Value = 0.7132873572634397
$ws.Range("B299").Value = 0
$ws.Range("B304").Value = 0
$ws.Range("B305").Value = 2.454299652530172
$ws.Range("B306").Value = 0
$ws.Range("B318").Value = 3.115238810587208
$ws.Range("B324").Value = 0.6463981650075452
$ws.Range("B325").Value = 0
$ws.Range("B334").Value = 0
$ws.Range("B336").Value = 0.5056632075671252
$ws.Range("B337").Value = 0
$ws.Range("B338").Value = 0.3600015234434402
$ws.Range("B340").Value = 0.392076332707383
$ws.Range("B341").Value = 0
$ws.Range("B342").Value = 0
$ws.Range("B343").Value = 0
$ws.Range("B349").Value = 0.001317988031008764
$ws.Range("B350").Value = 0.2894244325588191
$ws.Range("B354").Value = 0
$ws.Range("B356").Value = 0
$ws.Range("B358").Value = 0
$ws.Range("B359").Value = 1.449449116031491
